# Adds six new case rows (207-212) for case 21TRD09437 / Bunner, and
# removes the stray empty G205 cell, matching the updated Case_Data sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Helper: write a value that must stay plain text even when it looks
# numeric/currency-like (Excel would otherwise silently coerce "4510.11"
# or "$ 0" into a Number). Force the cell to Text format, assign, then
# strip the format back off so no stray style sticks around.
function Set-TextValue {
    param($cell, $value)
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.ClearFormats()
}

# Copying a whole A:K row pads in placeholder cells for every column in
# the copied range (even ones the source row never populated), so copy
# only the exact column spans each template row actually has - that
# keeps "missing" columns missing (no J/K here) instead of turning them
# into stray empty cells.

# 1) Seed row 210 from row 205 *before* row 205 is edited below - row
#    205 currently has an empty (but present) G cell, which is exactly
#    the shape row 210 needs for its own empty G210 cell. Row 205 has no
#    J/K worth keeping, so copy A:G and H:I only.
$ws.Range("A205:G205").Copy($ws.Range("A210"))
$ws.Range("H205:I205").Copy($ws.Range("H210"))

# 2) Seed rows 207, 209, 211, 212 from row 7, which already has exactly
#    the A:I layout (no J/K) that these rows need.
$ws.Range("A7:I7").Copy($ws.Range("A207"))
$ws.Range("A7:I7").Copy($ws.Range("A209"))
$ws.Range("A7:I7").Copy($ws.Range("A211"))
$ws.Range("A7:I7").Copy($ws.Range("A212"))

# 3) Seed row 208 from row 202, which has the A:F,H:I layout (no G, no
#    J/K) that row 208 needs - copy around the missing G so G208 is
#    never materialised.
$ws.Range("A202:F202").Copy($ws.Range("A208"))
$ws.Range("H202:I202").Copy($ws.Range("H208"))

# --- Row 207 ---
$ws.Range("A207").Value = "21TRD09437"
$ws.Range("B207").Value = "Bunner"
$ws.Range("C207").Value = "DUS"
Set-TextValue $ws.Range("D207") "4510.11"
$ws.Range("E207").Value = "M1"
$ws.Range("F207").Value = "No Contest"
$ws.Range("G207").Value = "Guilty"
Set-TextValue $ws.Range("H207") "$ 0"
Set-TextValue $ws.Range("I207") "$ 0"

# --- Row 208 ---
$ws.Range("A208").Value = "21TRD09437"
$ws.Range("B208").Value = "Bunner"
$ws.Range("C208").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
Set-TextValue $ws.Range("D208") "4511.21B1A"
$ws.Range("E208").Value = "M4"
$ws.Range("F208").Value = "Dismissed"
$ws.Range("H208").Value = " "
$ws.Range("I208").Value = " "

# --- Row 209 ---
$ws.Range("A209").Value = "21TRD09437"
$ws.Range("B209").Value = "Bunner"
$ws.Range("C209").Value = "RECKLESS OPERATION 1ST IN 1 YR"
Set-TextValue $ws.Range("D209") "4511.20"
$ws.Range("E209").Value = "MM"
$ws.Range("F209").Value = "No Contest"
$ws.Range("G209").Value = "Guilty"
Set-TextValue $ws.Range("H209") "$ 0"
Set-TextValue $ws.Range("I209") "$ 0"

# --- Row 210 ---
$ws.Range("A210").Value = "21TRD09437"
$ws.Range("B210").Value = "Bunner"
$ws.Range("C210").Value = "DUS"
Set-TextValue $ws.Range("D210") "4510.11"
$ws.Range("E210").Value = "M1"
$ws.Range("F210").Value = "Dismissed"
# G210 stays the empty cell copied from G205 above.
$ws.Range("H210").Value = " "
$ws.Range("I210").Value = " "

# --- Row 211 ---
$ws.Range("A211").Value = "21TRD09437"
$ws.Range("B211").Value = "Bunner"
$ws.Range("C211").Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
Set-TextValue $ws.Range("D211") "4511.21B1A"
$ws.Range("E211").Value = "M4"
$ws.Range("F211").Value = "No Contest"
$ws.Range("G211").Value = "Guilty"
Set-TextValue $ws.Range("H211") "$ 50"
Set-TextValue $ws.Range("I211") "$ 25"

# --- Row 212 ---
$ws.Range("A212").Value = "21TRD09437"
$ws.Range("B212").Value = "Bunner"
$ws.Range("C212").Value = "RECKLESS OPERATION 1ST IN 1 YR"
Set-TextValue $ws.Range("D212") "4511.20"
$ws.Range("E212").Value = "MM"
$ws.Range("F212").Value = "No Contest"
$ws.Range("G212").Value = "Guilty"
Set-TextValue $ws.Range("H212") "$ 0"
Set-TextValue $ws.Range("I212") "$ 0"

# 4) Drop the now-stray empty G205 cell (it disappears entirely in the
#    updated sheet, rather than remaining as an empty inline string).
$ws.Range("G205").Value = ""
